$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transition-probability matrix cells with recalculated values
# reflecting the larger simulated-game sample (per commit message: "added
# more games, sped up simulate game logic, and drafted optimization logic").
$ws.Range("B2").Value = 0.2182539682539683
$ws.Range("C2").Value = 0.5238095238095238
$ws.Range("J2").Value = 0.01984126984126984
$ws.Range("P2").Value = 0.1507936507936508
$ws.Range("S2").Value = 0.0873015873015873
$ws.Range("B3").Value = 0.007518796992481203
$ws.Range("C3").Value = 0.03007518796992481
$ws.Range("J3").Value = 0.03007518796992481
$ws.Range("P3").Value = 0.7218045112781954
$ws.Range("S3").Value = 0.2105263157894737
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.7419354838709677
$ws.Range("S4").Value = 0.1935483870967742
$ws.Range("B6").Value = 0.06521739130434782
$ws.Range("D6").Value = 0.02173913043478261
$ws.Range("F6").Value = 0.04347826086956522
$ws.Range("J6").Value = 0.2536231884057971
$ws.Range("O6").Value = 0.02173913043478261
$ws.Range("Q6").Value = 0.1884057971014493
$ws.Range("R6").Value = 0.07246376811594203
$ws.Range("B7").Value = 0.07228915662650602
$ws.Range("D7").Value = 0.006024096385542169
$ws.Range("F7").Value = 0.03614457831325301
$ws.Range("J7").Value = 0.1927710843373494
$ws.Range("O7").Value = 0.02409638554216868
$ws.Range("Q7").Value = 0.1686746987951807
$ws.Range("R7").Value = 0.0963855421686747
$ws.Range("S7").Value = 0.4036144578313253
$ws.Range("B8").Value = 0.0752212389380531
$ws.Range("D8").Value = 0.01548672566371681
$ws.Range("F8").Value = 0.0331858407079646
$ws.Range("J8").Value = 0.1128318584070796
$ws.Range("O8").Value = 0.01548672566371681
$ws.Range("Q8").Value = 0.2389380530973451
$ws.Range("R8").Value = 0.09734513274336283
$ws.Range("S8").Value = 0.411504424778761
$ws.Range("B9").Value = 0.0625
$ws.Range("D9").Value = 0.01785714285714286
$ws.Range("F9").Value = 0.01785714285714286
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.008928571428571428
$ws.Range("Q9").Value = 0.1964285714285714
$ws.Range("R9").Value = 0.07142857142857142
$ws.Range("S9").Value = 0.5
$ws.Range("B10").Value = 0.1048513302034429
$ws.Range("D10").Value = 0.01643192488262911
$ws.Range("F10").Value = 0.04773082942097027
$ws.Range("J10").Value = 0.136150234741784
$ws.Range("O10").Value = 0.01095461658841941
$ws.Range("Q10").Value = 0.2652582159624413
$ws.Range("R10").Value = 0.07902973395931143
$ws.Range("S10").Value = 0.3395931142410016
$ws.Range("G11").Value = 0.1580645161290322
$ws.Range("J11").Value = 0.1161290322580645
$ws.Range("K11").Value = 0.2354838709677419
$ws.Range("L11").Value = 0.4741935483870968
$ws.Range("S11").Value = 0.01612903225806452
$ws.Range("G12").Value = 0.6644736842105263
$ws.Range("J12").Value = 0.2960526315789473
$ws.Range("L12").Value = 0.01973684210526316
$ws.Range("S12").Value = 0.01973684210526316
$ws.Range("G13").Value = 0.4791666666666667
$ws.Range("J13").Value = 0.4375
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.01098901098901099
$ws.Range("H15").Value = 0.1923076923076923
$ws.Range("I15").Value = 0.03846153846153846
$ws.Range("J15").Value = 0.3846153846153846
$ws.Range("K15").Value = 0.04395604395604396
$ws.Range("M15").Value = 0.01098901098901099
$ws.Range("O15").Value = 0.05494505494505494
$ws.Range("S15").Value = 0.2637362637362637
$ws.Range("F16").Value = 0.01324503311258278
$ws.Range("H16").Value = 0.1456953642384106
$ws.Range("I16").Value = 0.06622516556291391
$ws.Range("J16").Value = 0.4304635761589404
$ws.Range("K16").Value = 0.1324503311258278
$ws.Range("M16").Value = 0.03973509933774835
$ws.Range("O16").Value = 0.02649006622516556
$ws.Range("S16").Value = 0.1456953642384106
$ws.Range("F17").Value = 0.01532567049808429
$ws.Range("H17").Value = 0.157088122605364
$ws.Range("I17").Value = 0.06896551724137931
$ws.Range("J17").Value = 0.446360153256705
$ws.Range("K17").Value = 0.09578544061302682
$ws.Range("M17").Value = 0.0210727969348659
$ws.Range("O17").Value = 0.05747126436781609
$ws.Range("S17").Value = 0.1379310344827586
$ws.Range("F18").Value = 0.01694915254237288
$ws.Range("H18").Value = 0.1242937853107345
$ws.Range("I18").Value = 0.0847457627118644
$ws.Range("J18").Value = 0.4350282485875706
$ws.Range("K18").Value = 0.1073446327683616
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.06214689265536723
$ws.Range("S18").Value = 0.1525423728813559
$ws.Range("F19").Value = 0.0113141862489121
$ws.Range("H19").Value = 0.2541340295909487
$ws.Range("I19").Value = 0.0391644908616188
$ws.Range("J19").Value = 0.3690165361183638
$ws.Range("K19").Value = 0.1174934725848564
$ws.Range("M19").Value = 0.02349869451697128
$ws.Range("O19").Value = 0.06353350739773717
$ws.Range("S19").Value = 0.1218450826805918
